# Insert a new data row at row 139 (pushing the existing rows 139-234 down to
# 140-235) and populate it with the new weekly price record, matching the
# commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 139; Excel shifts rows 139..234
# down to 140..235 and extends the used range to A1:R235 automatically.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new record.
$ws.Cells.Item(139, 1).Value = 4
$ws.Cells.Item(139, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(139, 3).Value = "Los Lagos"
$ws.Cells.Item(139, 4).Value = 44574
$ws.Cells.Item(139, 5).Value = 10
$ws.Cells.Item(139, 6).Value = 100114014
$ws.Cells.Item(139, 7).Value = "Betarraga"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 500
$ws.Cells.Item(139, 11).Value = 1000
$ws.Cells.Item(139, 12).Value = 1000
$ws.Cells.Item(139, 13).Value = 1000
$ws.Cells.Item(139, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(139, 15).Value = "Región del Maule"
$ws.Cells.Item(139, 16).Value = 200
$ws.Cells.Item(139, 17).Value = 5
$ws.Cells.Item(139, 18).Value = "Hortaliza"
